$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20: new date-entry row (copy date style from A19, then set the value)
$ws.Range("A19").Copy($ws.Range("A20"))
$ws.Range("A20").Value = 45862
$ws.Range("B20").Value = "Besprechung mit Gabriel wegen Entwürfe"
$ws.Range("G20").Value = 1

# Row 21: new date-entry row
$ws.Range("A19").Copy($ws.Range("A21"))
$ws.Range("A21").Value = 45863
$ws.Range("B21").Value = "Entwürfe überarbeitet + Filteransichten erstellt "
$ws.Range("G21").Value = 3

# Sheet view changes
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 123
$ws.Range("I22").Select()
